$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs right
# above the existing row 237 entry (chronologically it sits between the
# rows that are currently 236 and 237), so insert a fresh row at 237 -
# this pushes the old rows 237-287 down to 238-288 and grows the used
# range from A1:R287 to A1:R288.
$ws.Rows("237:237").Insert()

# Populate the newly inserted row 237 with the new record's data.
$ws.Range("A237").Value = 9
$ws.Range("B237").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C237").Value = "Metropolitana"
$ws.Range("D237").Value = 44511
$ws.Range("E237").Value = 13
$ws.Range("F237").Value = 100112031
$ws.Range("G237").Value = "Poroto verde"
$ws.Range("H237").Value = "Magnum"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 43
$ws.Range("K237").Value = 37000
$ws.Range("L237").Value = 40000
$ws.Range("M237").Value = 38465
$ws.Range("N237").Value = "$/malla 25 kilos"
$ws.Range("O237").Value = "Provincia de Limarí"
$ws.Range("P237").Value = 1539
$ws.Range("Q237").Value = 25
$ws.Range("R237").Value = "Hortaliza"
